# Timesheet.xlsx update - "Latest Commit on 20August from office"
#
# The sheet previously contained 5 daily status rows (rows 2-6). This
# edit replaces all of that with a single new day's entry (row 2) dated
# 08/16/2019, and removes the now-obsolete rows 3-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the four trailing status rows (old rows 3,4,5,6). Deleting row 3
# four times shifts everything up each time, ultimately leaving only the
# header row (1) and the single remaining data row (2).
$ws.Rows(3).EntireRow.Delete()
$ws.Rows(3).EntireRow.Delete()
$ws.Rows(3).EntireRow.Delete()
$ws.Rows(3).EntireRow.Delete()

# New remarks text for the remaining row, entered before the date so the
# shared-string table lists the long remark ahead of the date string.
$remarks = "1.`tVerified Branding page on Desktop and Mobile site. Testing is blocked for now due to some new pages showing under site map, updated just some time before, which have lots of issues. Please restore the pages to stable ones for further testing. `n2.`tIdentified some console errors in the application. Please find observations along with console errors in the document at: https://pmall4-my.sharepoint.com/:w:/g/personal/neerajv_pmall_com/ESEBSFeusalLj3EIWOfvezgB0imt1neXQp_iwQg7PB517Q?e=OhBA5z`n3.`tCreated 3 defects: #8015, #8016 and #8017 on dev environment.`n4.`tRetested defect #7732 on dev environment. Working fine on dev.`n5.`tUpdated status report for June release-cycle 2. Please find updated sheet attached."

$ws.Range("D2").Value = $remarks
$ws.Range("C2").Value = "08/16/2019"

# The new remark text needs a taller row than before.
$ws.Rows(2).RowHeight = 150
